$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '41.530.96'
$c.ClearFormats()
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  +0.35%  '
$c.ClearFormats()
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.482.25'
$c.ClearFormats()
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  +0.57%  '
$c.ClearFormats()
$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  -0.13%  '
$c.ClearFormats()
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '313.41'
$c.ClearFormats()
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  +0.50%  '
$c.ClearFormats()
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '92.44'
$c.ClearFormats()
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  -2.42%  '
$c.ClearFormats()
$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  -0.58%  '
$c.ClearFormats()
$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  -0.12%  '
$c.ClearFormats()
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.514'
$c.ClearFormats()
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  +2.76%  '
$c.ClearFormats()
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '32.91'
$c.ClearFormats()
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  -1.90%  '
$c.ClearFormats()
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  +1.42%  '
$c.ClearFormats()
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.111'
$c.ClearFormats()
$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  +2.06%  '
$c.ClearFormats()
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '2.863.92'
$c.ClearFormats()
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  +0.56%  '
$c.ClearFormats()
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '16.43'
$c.ClearFormats()
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  +9.67%  '
$c.ClearFormats()
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  -1.46%  '
$c.ClearFormats()
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '2.470.68'
$c.ClearFormats()
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  +3.01%  '
$c.ClearFormats()
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.776'
$c.ClearFormats()
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  -1.36%  '
$c.ClearFormats()
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '41.549.48'
$c.ClearFormats()
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  +0.36%  '
$c.ClearFormats()
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '6.55'
$c.ClearFormats()
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  +3.48%  '
$c.ClearFormats()
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '0.0₃0947'
$c.ClearFormats()
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  +2.66%  '
$c.ClearFormats()
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '72.54'
$c.ClearFormats()
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  +5.80%  '
$c.ClearFormats()
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '11.23'
$c.ClearFormats()
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  -0.11%  '
$c.ClearFormats()
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '237.11'
$c.ClearFormats()
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  +0.10%  '
$c.ClearFormats()
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '2.72'
$c.ClearFormats()
$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  -1.20%  '
$c.ClearFormats()
$c = $ws.Range('B25')
$c.NumberFormat = '@'
$c.Value = 'ImmutableX'
$c.ClearFormats()
$c = $ws.Range('C25')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c.ClearFormats()
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '1.91'
$c.ClearFormats()
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  +0.27%  '
$c.ClearFormats()
$c = $ws.Range('B26')
$c.NumberFormat = '@'
$c.Value = 'Dai'
$c.ClearFormats()
$c = $ws.Range('C26')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c.ClearFormats()
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.ClearFormats()
$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  -0.03%  '
$c.ClearFormats()
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '24.86'
$c.ClearFormats()
$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  +3.02%  '
$c.ClearFormats()
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '2.22'
$c.ClearFormats()
$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  +0.04%  '
$c.ClearFormats()
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '9.74'
$c.ClearFormats()
$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  +0.98%  '
$c.ClearFormats()
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '35.99'
$c.ClearFormats()
$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  -1.75%  '
$c.ClearFormats()
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '157.80'
$c.ClearFormats()
$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  +3.79%  '
$c.ClearFormats()
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '5.49'
$c.ClearFormats()
$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  -0.12%  '
$c.ClearFormats()
$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  -0.65%  '
$c.ClearFormats()
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.0759'
$c.ClearFormats()
$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  +1.44%  '
$c.ClearFormats()
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '17.46'
$c.ClearFormats()
$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  +2.34%  '
$c.ClearFormats()
$c = $ws.Range('B36')
$c.NumberFormat = '@'
$c.Value = 'ApeXProtocol'
$c.ClearFormats()
$c = $ws.Range('C36')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$c.ClearFormats()
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '2.38'
$c.ClearFormats()
$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  -10.21%  '
$c.ClearFormats()
$c = $ws.Range('B37')
$c.NumberFormat = '@'
$c.Value = 'Kaspa'
$c.ClearFormats()
$c = $ws.Range('C37')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c.ClearFormats()
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.106'
$c.ClearFormats()
$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  +3.95%  '
$c.ClearFormats()
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '2.92'
$c.ClearFormats()
$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  -4.33%  '
$c.ClearFormats()
$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  -2.71%  '
$c.ClearFormats()
$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  -0.01%  '
$c.ClearFormats()
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '4.09'
$c.ClearFormats()
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  -4.11%  '
$c.ClearFormats()
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  -0.27%  '
$c.ClearFormats()
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '1.973.93'
$c.ClearFormats()
$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  -0.73%  '
$c.ClearFormats()
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '19.15'
$c.ClearFormats()
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  -3.65%  '
$c.ClearFormats()
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.0285'
$c.ClearFormats()
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  -0.55%  '
$c.ClearFormats()
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  -2.43%  '
$c.ClearFormats()
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '8.94'
$c.ClearFormats()
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  +2.12%  '
$c.ClearFormats()
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '2.720.69'
$c.ClearFormats()
$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  +0.25%  '
$c.ClearFormats()
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '98.33'
$c.ClearFormats()
$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  +1.85%  '
$c.ClearFormats()
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '68.44'
$c.ClearFormats()
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  -2.03%  '
$c.ClearFormats()
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '72.62'
$c.ClearFormats()
$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  -3.00%  '
$c.ClearFormats()
